# Update Work Week and Social Spending
# -------------------------------------------------------------------------
# The "Data" sheet holds one row per (Country Code, Country Name,
# Indicator, Year) with the GDP-per-Capita figure in column E. This
# refreshes every existing year's figure (1950-2008) with the latest
# revision of the series, and appends the newly published years
# (2009-2016) as additional rows - each still stored as *text* (matching
# the source data, which keeps these numeric-looking readings as strings)
# rather than as a native number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$countryCode = 72
$countryName = "Botswana"
$indicator   = "GDP per Capita"
$firstYear   = 1950
$firstRow    = 2

# Year -> revised GDP-per-Capita reading (kept as text, same as source file)
$values = @(
    "489","497","504","513","521","529","537","544","550","559",
    "564","575","582","591","604","612","663","719","781","843",
    "907","1113","1361","1589","1656","1712","1801","1889","2005","2184",
    "2348","2528","2640","2848","3124","3301","3432","3601","3920","5066",
    "5267","5616.01123018046","5733.63002125501","5787.11450056963",
    "5946.42713209785","6318.47642374857","6647.7675526706",
    "7152.14743365942","7185.27225389974","7876.13851780229",
    "8044.14204336772","8082.82383350878","8604.06700130626",
    "9042.51973342875","9329.50921397854","9798.08706427939",
    "10665.6795265079","11602.0682465784","12385.1497252572",
    "11494.7565430293","12544.7273633297","13376","13743","14853",
    "15082","14805","15198"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $year = $firstYear + $i
    $row  = $firstRow + $i

    if ($year -gt 2008) {
        # Newly added year -> append a brand-new row.
        $ws.Cells.Item($row, 1).Value = $countryCode
        $ws.Cells.Item($row, 2).Value = $countryName
        $ws.Cells.Item($row, 3).Value = $indicator
        $ws.Cells.Item($row, 4).Value = $year
    }

    # Column E must stay text (these figures are numeric-looking strings
    # in the source, not real numbers) - force text formatting before the
    # assignment so Excel doesn't silently coerce it to a number, then
    # drop the formatting override again so no new cell style lingers.
    $cell = $ws.Cells.Item($row, 5)
    $cell.NumberFormat = "@"
    $cell.Value = $values[$i]
    $cell.ClearFormats()
}
